# Automatic map update (mapa_interactivo_INCO.html)
#
# The source feed dropped one record (Caso -173, "PACHECO DE MELO J A
# /ALT/ 2300") that existed between Caso -162 (row 14) and Caso -174
# (row 16). Removing that entire row shifts every subsequent record up
# by one and shrinks the used range from A1:P36 to A1:P35 - exactly
# matching the published diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(15).Delete()
